# Update from version 1.4 to 1.5:
# The test case "subjects" for TC1/TC3 are swapped, and TC2/TC4 are swapped.
# (TC5 - Perfis de Competencias - stays the same)
#
# TC1: Avaliacoes                    -> Competencias (portfolio)
# TC2: Niveis das Competencias       -> Periodos Avaliativos
# TC3: Competencias (portfolio)      -> Avaliacoes
# TC4: Periodos Avaliativos          -> Niveis das Competencias

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- TC1 block (rows 6-12): Avaliacoes -> Competencias (portfolio) ---
$ws.Range("B10").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Competencias (portfolio) a partir do menu inicial"
$ws.Range("D10").Value = "SYSTEM exibe a listagem das Competencias (portfolio) cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B11").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Competencias (portfolio)"

# --- TC2 block (rows 15-21): Niveis das Competencias -> Periodos Avaliativos ---
$ws.Range("B19").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Periodos Avaliativos a partir do menu inicial"
$ws.Range("D19").Value = "SYSTEM exibe a listagem dos Periodos Avaliativos cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B20").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Periodos Avaliativos"

# --- TC3 block (rows 24-30): Competencias (portfolio) -> Avaliacoes ---
$ws.Range("B28").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Avaliacoes a partir do menu inicial"
$ws.Range("D28").Value = "SYSTEM exibe a listagem das Avaliacoes cadastradas com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B29").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Avaliacoes"

# --- TC4 block (rows 33-39): Periodos Avaliativos -> Niveis das Competencias ---
$ws.Range("B37").Value = "Lider de Pessoas acessa a funcionalidade de gestao de Niveis das Competencias a partir do menu inicial"
$ws.Range("D37").Value = "SYSTEM exibe a listagem dos Niveis das Competencias cadastrados com opcoes de 'Novo', 'Editar', 'Excluir' e 'Ajuda'"
$ws.Range("B38").Value = "Lider de Pessoas clica na opcao 'Ajuda' para visualizar a ajuda da gestao de Niveis das Competencias"
